$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("chrY", 1043, 50, 0, "Pt0_blast", 0, "D1", "chrY_1043"),
    @("chrY", 1043, 50, 0, "Pt0_Er1", 0, "D1", "chrY_1043"),
    @("chrY", 1043, 50, 0, "Pt0_NK", 0, "D1", "chrY_1043"),
    @("chrY", 1043, 41, 9, "Pt0_Bcells", 0.18, "D1", "chrY_1043"),
    @("chr4", 502, 47, 0, "Pt0_blast", 0.06, "C1", "chr4_502"),
    @("chr4", 502, 48, 0, "Pt0_Er1", 0, "C1", "chr4_502"),
    @("chr4", 502, 44, 6, "Pt0_NK", 0.12, "C1", "chr4_502"),
    @("chr4", 502, 48, 0, "Pt0_Bcells", 0, "C1", "chr4_502"),
    @("chr9", 127, 43, 7, "Pt0_blast", 0.14, "B1", "chr9_127"),
    @("chr9", 127, 48, 2, "Pt0_Er1", 0.04, "B1", "chr9_127"),
    @("chr9", 127, 49, 0, "Pt0_NK", 0.02, "B1", "chr9_127"),
    @("chr9", 127, 49, 0, "Pt0_Bcells", 0, "B1", "chr9_127"),
    @("chr10", 500, 43, 7, "Pt0_blast", 0.14, "A1", "chr10_500"),
    @("chr10", 500, 34, 16, "Pt0_Er1", 0.32, "A1", "chr10_500"),
    @("chr10", 500, 45, 5, "Pt0_NK", 0.1, "A1", "chr10_500"),
    @("chr10", 500, 50, 0, "Pt0_Bcells", 0, "A1", "chr10_500")
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $r = $idx + 2
    $row = $data[$idx]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
}
